$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 167252
$ws.Range("C4").Value = 158147
$ws.Range("C5").Value = 9105
$ws.Range("C8").Value = 65.31999999999999
